# Update the "Förändrad" date column (C) from 45513 to 45514 for rows 2-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45513) {
        $cell.Value2 = 45514
    }
}
